$updates = @(
    @{ Cell = 'D2'; Value = '23.431.39' }
    @{ Cell = 'E2'; Value = '  -0.17%  ' }
    @{ Cell = 'D3'; Value = '1.631.33' }
    @{ Cell = 'E3'; Value = '  -0.44%  ' }
    @{ Cell = 'D4'; Value = '1.001' }
    @{ Cell = 'E4'; Value = '  -0.10%  ' }
    @{ Cell = 'E5'; Value = '  -0.02%  ' }
    @{ Cell = 'D6'; Value = '304.74' }
    @{ Cell = 'E6'; Value = '  -0.77%  ' }
    @{ Cell = 'D7'; Value = '0.3781' }
    @{ Cell = 'E7'; Value = '  +0.24%  ' }
    @{ Cell = 'D8'; Value = '0.3654' }
    @{ Cell = 'E8'; Value = '  +0.21%  ' }
    @{ Cell = 'D9'; Value = '51.59' }
    @{ Cell = 'E9'; Value = '  -1.25%  ' }
    @{ Cell = 'D10'; Value = '0.08248' }
    @{ Cell = 'E10'; Value = '  +0.92%  ' }
    @{ Cell = 'D11'; Value = '1.232' }
    @{ Cell = 'E11'; Value = '  -2.79%  ' }
    @{ Cell = 'E12'; Value = '  -0.09%  ' }
    @{ Cell = 'D13'; Value = '22.42' }
    @{ Cell = 'E13'; Value = '  -2.54%  ' }
    @{ Cell = 'D14'; Value = '6.563' }
    @{ Cell = 'E14'; Value = '  -1.23%  ' }
    @{ Cell = 'D15'; Value = '0.00001254' }
    @{ Cell = 'E15'; Value = '  -1.98%  ' }
    @{ Cell = 'D16'; Value = '7.343' }
    @{ Cell = 'E16'; Value = '  -0.54%  ' }
    @{ Cell = 'D17'; Value = '1.628.47' }
    @{ Cell = 'E17'; Value = '  -0.86%  ' }
    @{ Cell = 'D18'; Value = '94.13' }
    @{ Cell = 'E18'; Value = '  -0.78%  ' }
    @{ Cell = 'D19'; Value = '0.06981' }
    @{ Cell = 'E19'; Value = '  +0.22%  ' }
    @{ Cell = 'E20'; Value = '  -2.63%  ' }
    @{ Cell = 'D21'; Value = '6.534' }
    @{ Cell = 'E21'; Value = '  -0.28%  ' }
    @{ Cell = 'E22'; Value = '  +0.02%  ' }
    @{ Cell = 'D23'; Value = '12.75' }
    @{ Cell = 'E23'; Value = '  -0.53%  ' }
    @{ Cell = 'D24'; Value = '23.420.81' }
    @{ Cell = 'E24'; Value = '  -0.20%  ' }
    @{ Cell = 'D25'; Value = '3.206' }
    @{ Cell = 'E25'; Value = '  +3.14%  ' }
    @{ Cell = 'D26'; Value = '2.470' }
    @{ Cell = 'E26'; Value = '  +1.96%  ' }
    @{ Cell = 'D27'; Value = '21.43' }
    @{ Cell = 'E27'; Value = '  +0.75%  ' }
    @{ Cell = 'D28'; Value = '150.61' }
    @{ Cell = 'E28'; Value = '  -0.52%  ' }
    @{ Cell = 'D29'; Value = '5.299' }
    @{ Cell = 'E29'; Value = '  -1.24%  ' }
    @{ Cell = 'D30'; Value = '134.39' }
    @{ Cell = 'E30'; Value = '  -0.95%  ' }
    @{ Cell = 'D31'; Value = '1.809.81' }
    @{ Cell = 'E31'; Value = '  -0.80%  ' }
    @{ Cell = 'D32'; Value = '2.262' }
    @{ Cell = 'E32'; Value = '  -3.69%  ' }
    @{ Cell = 'D33'; Value = '6.835' }
    @{ Cell = 'E33'; Value = '  +0.71%  ' }
    @{ Cell = 'D34'; Value = '1.021' }
    @{ Cell = 'E34'; Value = '  +5.78%  ' }
    @{ Cell = 'D35'; Value = '10.92' }
    @{ Cell = 'E35'; Value = '  +5.75%  ' }
    @{ Cell = 'D36'; Value = '0.02798' }
    @{ Cell = 'E36'; Value = '  -1.09%  ' }
    @{ Cell = 'D37'; Value = '0.2531' }
    @{ Cell = 'E37'; Value = '  -0.26%  ' }
    @{ Cell = 'D38'; Value = '0.08788' }
    @{ Cell = 'E38'; Value = '  -0.79%  ' }
    @{ Cell = 'B39'; Value = 'InternetComputer(DFINITY)' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Cell = 'D39'; Value = '6.088' }
    @{ Cell = 'E39'; Value = '  -1.43%  ' }
    @{ Cell = 'B40'; Value = 'Hedera' }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = 'D40'; Value = '0.07138' }
    @{ Cell = 'E40'; Value = '  -2.97%  ' }
    @{ Cell = 'D41'; Value = '1.354' }
    @{ Cell = 'E41'; Value = '  -2.28%  ' }
    @{ Cell = 'D42'; Value = '0.7051' }
    @{ Cell = 'E42'; Value = '  -0.84%  ' }
    @{ Cell = 'D43'; Value = '16.08' }
    @{ Cell = 'E43'; Value = '  -1.19%  ' }
    @{ Cell = 'D44'; Value = '12.22' }
    @{ Cell = 'E44'; Value = '  -2.39%  ' }
    @{ Cell = 'D45'; Value = '0.6577' }
    @{ Cell = 'E45'; Value = '  +0.29%  ' }
    @{ Cell = 'D46'; Value = '2.321' }
    @{ Cell = 'E46'; Value = '  -0.84%  ' }
    @{ Cell = 'E47'; Value = '  -0.03%  ' }
    @{ Cell = 'D48'; Value = '3.983' }
    @{ Cell = 'E48'; Value = '  -1.13%  ' }
    @{ Cell = 'D49'; Value = '0.08014' }
    @{ Cell = 'E49'; Value = '  +0.61%  ' }
    @{ Cell = 'D50'; Value = '1.200' }
    @{ Cell = 'E50'; Value = '  -0.75%  ' }
    @{ Cell = 'D51'; Value = '126.48' }
    @{ Cell = 'E51'; Value = '  -2.28%  ' }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($u in $updates) {
    $cellRef = $u.Cell
    $val = $u.Value
    $col = $cellRef -replace '[0-9]+$', ''
    $rng = $ws.Range($cellRef)

    if ($col -eq "D") {
        # Column D values are plain-text numeric-looking strings (e.g. "1.001",
        # "23.431.39"). Force text storage so Excel doesn't reinterpret them as
        # numbers (which would strip trailing zeros / use two-dot forms, etc).
        $rng.NumberFormat = "@"
        $rng.Value = $val
        $rng.Style = "Normal"
    } else {
        $rng.Value = $val
    }
}
